# Revision previa a envio SETI produccion 3
# Se actualiza la columna de "NEW" que indica a SETI cambios en posiciones

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BATT_CPU")
$ws2 = $wb.Worksheets.Item("_HISTORY")

# --- BATT_CPU: flip the "NEW" (YES/NO) flags for the positions that moved
#     relative to the previous SETI revision (7) ---
$ws1.Range("B8").Value  = "YES"   # position 2  (U? part) now flagged as new
$ws1.Range("B12").Value = "NO"    # position 6 no longer flagged as new
$ws1.Range("B21").Value = "YES"   # position 15 now flagged as new
$ws1.Range("B37").Value = "YES"   # position 31 now flagged as new

# --- _HISTORY: log this revision as a new row, copying the format of the
#     previous row (date format on column B, wrapped text on column D) ---
$ws2.Range("B12").Copy($ws2.Range("B13"))
$ws2.Range("D12").Copy($ws2.Range("D13"))

$ws2.Cells.Item(13, 2).Value = "9/26/2023"
$ws2.Cells.Item(13, 3).Value = "DGB"
$ws2.Cells.Item(13, 4).Value = "Se actualiza la columna de NEW para que SETI tenga los cambios actualizados contra la version anterior (7)"

# --- restore the on-screen selections left by the editor ---
[void]$ws2.Range("D14").Select()
[void]$ws1.Range("C3").Select()
